# Update viewer counts (F) and minimum ticket prices (G) to refreshed
# values as scraped at commit 456a3b4, across all four sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 14593
$ws.Range("G3").Value = 80
$ws.Range("F4").Value = 14818
$ws.Range("G4").Value = 109
$ws.Range("G5").Value = 188
$ws.Range("F7").Value = 5997
$ws.Range("G7").Value = 68
$ws.Range("F13").Value = 1596
$ws.Range("F14").Value = 475
$ws.Range("F16").Value = 1265
$ws.Range("F17").Value = 1920
$ws.Range("F18").Value = 933
$ws.Range("F20").Value = 2325
$ws.Range("F23").Value = 3489
$ws.Range("F24").Value = 312
$ws.Range("F25").Value = 332
$ws.Range("F26").Value = 2563
$ws.Range("F27").Value = 641
$ws.Range("F30").Value = 1871
$ws.Range("F31").Value = 1112
$ws.Range("F32").Value = 1525
$ws.Range("F35").Value = 7219
$ws.Range("F36").Value = 5059
$ws.Range("F38").Value = 706
$ws.Range("F39").Value = 710
$ws.Range("F40").Value = 3352
$ws.Range("F44").Value = 140
$ws.Range("F45").Value = 109
$ws.Range("F46").Value = 4477
$ws.Range("F47").Value = 665
$ws.Range("F48").Value = 11
$ws.Range("F49").Value = 320

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 136
$ws.Range("F22").Value = 75
$ws.Range("F29").Value = 17

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 7853
$ws.Range("F3").Value = 292
$ws.Range("F4").Value = 1015

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 7853
$ws.Range("F4").Value = 292
$ws.Range("F5").Value = 1015
$ws.Range("F7").Value = 14593
$ws.Range("G7").Value = 80
$ws.Range("F8").Value = 14818
$ws.Range("G8").Value = 109
$ws.Range("G10").Value = 188
$ws.Range("F12").Value = 5997
$ws.Range("G12").Value = 68
$ws.Range("F13").Value = 136
$ws.Range("F17").Value = 1596
$ws.Range("F18").Value = 475
$ws.Range("F19").Value = 1265
$ws.Range("F22").Value = 3489
$ws.Range("F23").Value = 332
$ws.Range("F24").Value = 2563
$ws.Range("F25").Value = 641
$ws.Range("F27").Value = 1871
$ws.Range("F33").Value = 1112
$ws.Range("F35").Value = 7219
$ws.Range("F36").Value = 5059
$ws.Range("F38").Value = 706
$ws.Range("F39").Value = 3352
$ws.Range("F42").Value = 140
$ws.Range("F44").Value = 109
$ws.Range("F45").Value = 665
$ws.Range("F46").Value = 320

Write-Output "Updated $(($wb.Worksheets | Measure-Object).Count) sheets with refreshed attendance/price figures."
